$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "202.173.124.249"
$ws.Range("B6").Value = 28.5212672
$ws.Range("C6").Value = 77.2243456
$ws.Range("D6").Value = 735485.4806669627
$ws.Range("E6").Value = "Mozilla/5.0 (Windows NT 10.0; Win64; x64) AppleWebKit/537.36 (KHTML, like Gecko) Chrome/137.0.0.0 Safari/537.36"
$ws.Range("F6").Value = "Win32"
$ws.Range("G6").Value = "2025-06-21T04:23:56.793Z"

$ws.Range("A7").Value = "106.219.230.230"
$ws.Range("B7").Value = 28.3621629
$ws.Range("C7").Value = 77.2827333
$ws.Range("D7").Value = 12.9399995803833
$ws.Range("E7").Value = "Mozilla/5.0 (Linux; Android 10; K) AppleWebKit/537.36 (KHTML, like Gecko) Chrome/137.0.0.0 Mobile Safari/537.36"
$ws.Range("F7").Value = "Linux armv81"
$ws.Range("G7").Value = "2025-06-21T04:24:08.731Z"

$ws.Range("A8").Value = "106.219.230.230"
$ws.Range("B8").Value = 28.3621629
$ws.Range("C8").Value = 77.2827333
$ws.Range("D8").Value = 55.656795501708984
$ws.Range("E8").Value = "Mozilla/5.0 (Linux; Android 10; K) AppleWebKit/537.36 (KHTML, like Gecko) Chrome/137.0.0.0 Mobile Safari/537.36"
$ws.Range("F8").Value = "Linux armv81"
$ws.Range("G8").Value = "2025-06-21T04:24:37.801Z"
